$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the weekly data block (rows 411-412),
# shifting all the existing data (rows 411-536) down by two rows.
$ws.Rows("411:412").Insert()

# Populate the two newly inserted rows with the new week's data.

# Row 411 - "Primera" quality
$ws.Cells.Item(411, 1).Value2  = 8
$ws.Cells.Item(411, 2).Value2  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(411, 3).Value2  = "Coquimbo"
$ws.Cells.Item(411, 4).Value2  = 44900
$ws.Cells.Item(411, 5).Value2  = 4
$ws.Cells.Item(411, 6).Value2  = 100112009
$ws.Cells.Item(411, 7).Value2  = "Acelga"
$ws.Cells.Item(411, 8).Value2  = "Sin especificar"
$ws.Cells.Item(411, 9).Value2  = "Primera"
$ws.Cells.Item(411, 10).Value2 = 2200
$ws.Cells.Item(411, 11).Value2 = 650
$ws.Cells.Item(411, 12).Value2 = 700
$ws.Cells.Item(411, 13).Value2 = 675
$ws.Cells.Item(411, 14).Value2 = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(411, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(411, 16).Value2 = 338
$ws.Cells.Item(411, 17).Value2 = 2
$ws.Cells.Item(411, 18).Value2 = "Hortaliza"

# Row 412 - "Segunda" quality
$ws.Cells.Item(412, 1).Value2  = 8
$ws.Cells.Item(412, 2).Value2  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(412, 3).Value2  = "Coquimbo"
$ws.Cells.Item(412, 4).Value2  = 44900
$ws.Cells.Item(412, 5).Value2  = 4
$ws.Cells.Item(412, 6).Value2  = 100112009
$ws.Cells.Item(412, 7).Value2  = "Acelga"
$ws.Cells.Item(412, 8).Value2  = "Sin especificar"
$ws.Cells.Item(412, 9).Value2  = "Segunda"
$ws.Cells.Item(412, 10).Value2 = 1460
$ws.Cells.Item(412, 11).Value2 = 550
$ws.Cells.Item(412, 12).Value2 = 600
$ws.Cells.Item(412, 13).Value2 = 575
$ws.Cells.Item(412, 14).Value2 = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(412, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(412, 16).Value2 = 288
$ws.Cells.Item(412, 17).Value2 = 2
$ws.Cells.Item(412, 18).Value2 = "Hortaliza"
